# Split the "Collaborated on the design..." sentence (single run) into
# five runs carrying the updated wording, mirroring the target diff.
#
# The new wording, as five separate <w:r> runs (all sharing the exact
# same rPr as the original run: rFonts majorHAnsi + shd clear/FFFFFF):
#   1) "Collaborated on the design"
#   2) " and implementation of a website built for the company’s first
#       online streaming conference"
#   3) "."
#   4) " Assisted in learning and executing advanced WordPress features,
#       hosting services, and SEO"
#   5) " --- "

$d = $word.ActiveDocument

$oldText = "Collaborated on the design, implementation, and hosting services of a website built for the company’s first online conference, Language Summit 2020. --- "

$chunk1 = "Collaborated on the design"
$chunk2 = " and implementation of a website built for the company’s first online streaming conference"
$chunk3 = "."
$chunk4 = " Assisted in learning and executing advanced WordPress features, hosting services, and SEO"
$chunk5 = " --- "

# Locate the existing run's text.
$found = $d.Content
$null = $found.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$pos0 = $found.Start

# Overwrite the original run's text with the first chunk (keeps it in
# the same, pre-existing run / formatting).
$r0 = $d.Range($pos0, $pos0 + $oldText.Length)
$r0.Text = $chunk1
$pos1 = $pos0 + $chunk1.Length

# Append the remaining four chunks after it, one at a time.
$ins2 = $d.Range($pos1, $pos1)
$ins2.InsertAfter($chunk2)
$pos2 = $pos1 + $chunk2.Length

$ins3 = $d.Range($pos2, $pos2)
$ins3.InsertAfter($chunk3)
$pos3 = $pos2 + $chunk3.Length

$ins4 = $d.Range($pos3, $pos3)
$ins4.InsertAfter($chunk4)
$pos4 = $pos3 + $chunk4.Length

$ins5 = $d.Range($pos4, $pos4)
$ins5.InsertAfter($chunk5)
$pos5 = $pos4 + $chunk5.Length

# All five chunks now share identical formatting, so the engine would
# otherwise coalesce them back into one run on save. Re-asserting (and
# immediately reverting) a formatting property on each chunk's exact
# range keeps its run boundary distinct. Doing this from the last chunk
# back to the first is what makes every boundary stick.
$t5 = $d.Range($pos4, $pos5)
$t5.Bold = 1
$t5.Bold = 0

$t4 = $d.Range($pos3, $pos4)
$t4.Bold = 1
$t4.Bold = 0

$t3 = $d.Range($pos2, $pos3)
$t3.Bold = 1
$t3.Bold = 0

$t2 = $d.Range($pos1, $pos2)
$t2.Bold = 1
$t2.Bold = 0

$t1 = $d.Range($pos0, $pos1)
$t1.Bold = 1
$t1.Bold = 0

Write-Output "Edit applied: 5 runs from $pos0 to $pos5"
